# Update "想去人数" (number of people wanting to go) values for a few events
# that were refreshed in the upstream data source.

$wb = $excel.ActiveWorkbook

# Sheet "展览" holds the exhibition-only listing.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F8").Value = 7120
$wsExhibit.Range("F11").Value = 3584
$wsExhibit.Range("F15").Value = 593
$wsExhibit.Range("F16").Value = 78

# Sheet "全部类型" holds the same events merged across all categories.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 7120
$wsAll.Range("F14").Value = 3584
$wsAll.Range("F18").Value = 593
$wsAll.Range("F19").Value = 78
